# ajustes fin de mes enero
# - Remove the row for placa "33-88" (Leidy) from every monthly sheet
#   (ene2025..jun2025); all rows below it shift up by one.
# - On ene2025, mark every resident's January dues as fully paid
#   (cuota column C copied into the "pago2" column D), except the
#   reduced-quota resident (now row 12) who pays 56000.
# - Update each sheet's remembered cell selection.
# - Tidy the "pagos" sheet header/footer font style name.

$wb = $excel.ActiveWorkbook

# --- Monthly sheets: delete the "33-88 / Leidy" row (row 12) ---
$monthSheets = @("ene2025", "feb2025", "mar2025", "abr2025", "may2025", "jun2025")
foreach ($name in $monthSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Rows.Item(12).Delete()
}

# --- ene2025: close out January, everyone paid in full ---
$wsEne = $wb.Worksheets.Item("ene2025")
$wsEne.Range("C2:C22").Value = 65000
$wsEne.Range("D2:D22").Value = 65000
$wsEne.Range("C12").Value = 56000
$wsEne.Range("D12").Value = 56000

# --- Selection bookkeeping per sheet ---
$wb.Worksheets.Item("feb2025").Range("A12").Select() | Out-Null
$wb.Worksheets.Item("mar2025").Range("A12").Select() | Out-Null
$wb.Worksheets.Item("abr2025").Range("A12").Select() | Out-Null
$wb.Worksheets.Item("may2025").Range("A12").Select() | Out-Null
$wb.Worksheets.Item("jun2025").Range("F25").Select() | Out-Null

# ene2025 is the active tab, so select it last.
$wsEne.Range("E28").Select() | Out-Null

# --- pagos: header/footer font style "Regular" -> "Normal" ---
$wsPagos = $wb.Worksheets.Item("pagos")
$wsPagos.PageSetup.CenterHeader = '&"Times New Roman,Normal"&12&A'
$wsPagos.PageSetup.CenterFooter = '&"Times New Roman,Normal"&12Página &P'
